$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.339005708694458
$ws.Range("B1").Value = 1.252324461936951
$ws.Range("C1").Value = 5.401671886444092
$ws.Range("D1").Value = 2.167785882949829
$ws.Range("E1").Value = 1.130434393882751
